# basic_DMI.xlsx - write DMI (Directional Movement Index) headers to Excel
#
# The sheet originally had the 最高價(High)/最低價(Low)/收盤價(Close) headers in
# B1:D1 out of their natural order (收盤價, 最高價, 最低價) and the last two
# header columns (S1, T1) were left blank even though the workbook already
# dimensions the sheet through column T. This finishes the DMI header row by
# putting High/Low/Close back in their natural left-to-right order and adding
# the missing "DXt" / "ADXt" columns, pushing the trailing
# "DMO(Directional Movement Oscillator)" label out to column T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the price headers: High, Low, Close
$ws.Range("B1").Value = "最高價"
$ws.Range("C1").Value = "最低價"
$ws.Range("D1").Value = "收盤價"

# DXt / ADXt are new; the DMO label moves from R1 out to T1
$ws.Range("R1").Value = "DXt"
$ws.Range("S1").Value = "ADXt"
$ws.Range("T1").Value = "DMO(Directional Movement Oscillator)"

# Column T now carries the long DMO label, so fit it to the new content
$ws.Columns("T").AutoFit()

# Leave the selection where the author ended up when finishing the row
[void]$ws.Range("G14").Select()
